$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared text strings (volume number and report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/18/2023  Through  9/24/2023"

# --- Update weekly crime-statistics figures (rows 14-30) ---
# Row 14
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 8
$ws.Range("H14").Value = -37.5
$ws.Range("I14").Value = 51
$ws.Range("J14").Value = 54
$ws.Range("K14").Value = -5.555555555555
$ws.Range("L14").Value = -26.086956521739
$ws.Range("M14").Value = -49.504950495049
$ws.Range("N14").Value = -86.027397260274

# Row 15
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 14
$ws.Range("G15").Value = 28
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 163
$ws.Range("J15").Value = 188
$ws.Range("K15").Value = -13.297872340425
$ws.Range("L15").Value = 0.617283950617
$ws.Range("M15").Value = 7.236842105263
$ws.Range("N15").Value = -64.175824175824

# Row 16
$ws.Range("C16").Value = 56
$ws.Range("D16").Value = 51
$ws.Range("E16").Value = 9.803921568627
$ws.Range("F16").Value = 200
$ws.Range("G16").Value = 222
$ws.Range("H16").Value = -9.909909909909
$ws.Range("I16").Value = 1791
$ws.Range("J16").Value = 1918
$ws.Range("K16").Value = -6.621480709071
$ws.Range("L16").Value = 23.262216104611
$ws.Range("M16").Value = -28.758949880668
$ws.Range("N16").Value = -85.062552126772

# Row 17
$ws.Range("D17").Value = 87
$ws.Range("E17").Value = -9.195402298850
$ws.Range("F17").Value = 326
$ws.Range("G17").Value = 345
$ws.Range("H17").Value = -5.507246376811
$ws.Range("I17").Value = 3125
$ws.Range("J17").Value = 3089
$ws.Range("K17").Value = 1.165425704111
$ws.Range("L17").Value = 20.843000773395
$ws.Range("M17").Value = 25.401284109149
$ws.Range("N17").Value = -50.903377847604

# Row 18
$ws.Range("C18").Value = 44
$ws.Range("D18").Value = 55
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 150
$ws.Range("G18").Value = 171
$ws.Range("H18").Value = -12.280701754386
$ws.Range("I18").Value = 1504
$ws.Range("J18").Value = 1753
$ws.Range("K18").Value = -14.204221334854
$ws.Range("L18").Value = 5.027932960893
$ws.Range("M18").Value = -34.551784160139
$ws.Range("N18").Value = -83.068783068783

# Row 19
$ws.Range("C19").Value = 115
$ws.Range("D19").Value = 118
$ws.Range("E19").Value = -2.542372881355
$ws.Range("F19").Value = 444
$ws.Range("G19").Value = 493
$ws.Range("H19").Value = -9.939148073022
$ws.Range("I19").Value = 4243
$ws.Range("J19").Value = 4344
$ws.Range("K19").Value = -2.325046040515
$ws.Range("L19").Value = 29.676039119804
$ws.Range("M19").Value = 38.796205430160
$ws.Range("N19").Value = -15.796785076404

# Row 20
$ws.Range("C20").Value = 40
$ws.Range("D20").Value = 41
$ws.Range("E20").Value = -2.439024390243
$ws.Range("F20").Value = 153
$ws.Range("G20").Value = 168
$ws.Range("H20").Value = -8.928571428571
$ws.Range("I20").Value = 1340
$ws.Range("J20").Value = 1338
$ws.Range("K20").Value = 0.149476831091
$ws.Range("L20").Value = 22.374429223744
$ws.Range("M20").Value = 28.846153846153
$ws.Range("N20").Value = -80.455075845974

# Row 21
$ws.Range("C21").Value = 338
$ws.Range("D21").Value = 362
$ws.Range("E21").Value = -6.629834254143
$ws.Range("F21").Value = 1292
$ws.Range("G21").Value = 1435
$ws.Range("H21").Value = -9.965156794425
$ws.Range("I21").Value = 12217
$ws.Range("J21").Value = 12684
$ws.Range("K21").Value = -3.681803847366
$ws.Range("L21").Value = 21.332803654782
$ws.Range("M21").Value = 4.830959327269
$ws.Range("N21").Value = -69.421570345155

# Row 22
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 23
$ws.Range("G22").Value = 26
$ws.Range("H22").Value = -11.538461538461
$ws.Range("I22").Value = 212
$ws.Range("J22").Value = 257
$ws.Range("K22").Value = -17.509727626459
$ws.Range("L22").Value = 19.774011299435
$ws.Range("M22").Value = -30.263157894736

# Row 23
$ws.Range("C23").Value = 31
$ws.Range("D23").Value = 35
$ws.Range("E23").Value = -11.428571428571
$ws.Range("F23").Value = 108
$ws.Range("G23").Value = 124
$ws.Range("H23").Value = -12.903225806451
$ws.Range("I23").Value = 1158
$ws.Range("J23").Value = 1138
$ws.Range("K23").Value = 1.757469244288
$ws.Range("L23").Value = 9.659090909090
$ws.Range("M23").Value = 34.338747099768

# Row 24
$ws.Range("C24").Value = 287
$ws.Range("D24").Value = 345
$ws.Range("E24").Value = -16.811594202898
$ws.Range("F24").Value = 1029
$ws.Range("G24").Value = 1151
$ws.Range("H24").Value = -10.599478714161
$ws.Range("I24").Value = 9229
$ws.Range("J24").Value = 9832
$ws.Range("K24").Value = -6.133034987794
$ws.Range("L24").Value = 26.338124572210
$ws.Range("M24").Value = 22.954969357847

# Row 25
$ws.Range("C25").Value = 114
$ws.Range("D25").Value = 124
$ws.Range("E25").Value = -8.064516129032
$ws.Range("F25").Value = 509
$ws.Range("G25").Value = 455
$ws.Range("H25").Value = 11.868131868131
$ws.Range("I25").Value = 4518
$ws.Range("J25").Value = 4401
$ws.Range("K25").Value = 2.658486707566
$ws.Range("L25").Value = 34.224598930481
$ws.Range("M25").Value = -23.617920540997

# Row 26
$ws.Range("C26").Value = 4
$ws.Range("E26").Value = -63.636363636363
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = -40
$ws.Range("I26").Value = 254
$ws.Range("J26").Value = 281
$ws.Range("K26").Value = -9.608540925266
$ws.Range("L26").Value = -9.608540925266

# Row 27
$ws.Range("C27").Value = 19
$ws.Range("D27").Value = 7
$ws.Range("E27").Value = 171.428571428571
$ws.Range("G27").Value = 40
$ws.Range("H27").Value = 60
$ws.Range("I27").Value = 477
$ws.Range("J27").Value = 453
$ws.Range("K27").Value = 5.298013245033
$ws.Range("L27").Value = -5.357142857142

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = -66.666666666666
$ws.Range("F28").Value = 21
$ws.Range("G28").Value = 27
$ws.Range("H28").Value = -22.222222222222
$ws.Range("I28").Value = 175
$ws.Range("J28").Value = 262
$ws.Range("K28").Value = -33.206106870229
$ws.Range("L28").Value = -45.987654320987
$ws.Range("M28").Value = -56.575682382134
$ws.Range("N28").Value = -87.897648686030

# Row 29
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = -60
$ws.Range("F29").Value = 14
$ws.Range("G29").Value = 20
$ws.Range("H29").Value = -30
$ws.Range("I29").Value = 149
$ws.Range("J29").Value = 219
$ws.Range("K29").Value = -31.963470319634
$ws.Range("L29").Value = -42.023346303501
$ws.Range("M29").Value = -54.012345679012
$ws.Range("N29").Value = -88.547271329746

# Row 30
$ws.Range("D30").Value = 4
$ws.Range("J30").Value = 61
$ws.Range("K30").Value = -24.590163934426
$ws.Range("L30").Value = -6.122448979591

